# Generate Report for Handback
# - Updates the "Status" shared string so every cell referencing it
#   (Overview!E/F + zh-cn!C + de-de!C) flips from "Ready for handoff" to
#   "Handed back: in sync with en-US".
# - Fills in "Latest Target File" (col I) and "Latest Handback File" (col J)
#   on the zh-cn / de-de detail sheets, with a hyperlink on the new I cell.
# - Stamps the "Latest Handback DateTime" (col K) with the handback time
#   (same value on zh-cn, a later value on de-de).
# - Widens a couple of columns that now hold longer text.

$wb = $excel.ActiveWorkbook

$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bc1b4b95be31d17edd296968c280b7a255a541e7/e2e/a.md"

# ---- Status text, shared across Overview/zh-cn/de-de ----------------------
$wb.Worksheets.Item("zh-cn").Range("C2").Value = "Handed back: in sync with en-US"
$wb.Worksheets.Item("zh-cn").Range("C3").Value = "Handed back: in sync with en-US"
$wb.Worksheets.Item("de-de").Range("C2").Value = "Handed back: in sync with en-US"
$wb.Worksheets.Item("de-de").Range("C3").Value = "Handed back: in sync with en-US"

# ---- zh-cn detail sheet -----------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("I2").Value = "a.md"
$ws.Hyperlinks.Add($ws.Range("I2"), $targetUrl, "", "", "a.md")
$ws.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$ws.Range("K2").Value = "2016-08-27 08:36:43"

$ws.Range("I3").Value = "a.md"
$ws.Hyperlinks.Add($ws.Range("I3"), $targetUrl, "", "", "a.md")
$ws.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$ws.Range("K3").Value = "2016-08-27 08:36:43"

$ws.Columns.Item(3).ColumnWidth = 29.1
$ws.Columns.Item(10).ColumnWidth = 39.2

# ---- de-de detail sheet -----------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("I2").Value = "a.md"
$ws.Hyperlinks.Add($ws.Range("I2"), $targetUrl, "", "", "a.md")
$ws.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$ws.Range("K2").Value = "2016-08-27 08:36:50"

$ws.Range("I3").Value = "a.md"
$ws.Hyperlinks.Add($ws.Range("I3"), $targetUrl, "", "", "a.md")
$ws.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$ws.Range("K3").Value = "2016-08-27 08:36:50"

$ws.Columns.Item(3).ColumnWidth = 29.1
$ws.Columns.Item(10).ColumnWidth = 39.2

# ---- Overview sheet: widen the columns that mirror the Status text -----
$ws = $wb.Worksheets.Item("Overview")
$ws.Columns.Item(5).ColumnWidth = 29.1
$ws.Columns.Item(6).ColumnWidth = 29.1
